# Hortaliza, Femacal de La Calera - Perejil
# Weekly update: insert a new, most-recent observation as row 7
# (pushing the previous rows 7..37 down to 8..38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7; existing rows 7-37 shift down to 8-38.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly record.
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "Femacal de La Calera"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 45189
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 100112044
$ws.Range("G7").Value = "Perejil"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 3000
$ws.Range("N7").Value = "$/docena de atados (3 kilos)"
$ws.Range("O7").Value = "Provincia de Quillota"
$ws.Range("P7").Value = 1000
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = "Hortaliza"
